$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Cells.Item(2, 8)
Write-Host $cell.Value()
Write-Host $cell.Formula()
